# Append a new log entry (row 47) to the time-tracking sheet:
#   Date = 2025-01-29 (serial 45686), Tasks done = "city art", hours = 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = 45686
$ws.Range("A47").NumberFormat = "d-mmm"
$ws.Range("B47").Value = "city art"
$ws.Range("C47").Value = 4

# Match the author's final selection, which moved to the newly typed cell.
$ws.Range("B47").Select()
